$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "ahmed"
$ws.Range("B5").Value = $true
